# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (totals) sheet
#    and populate it with the Q1-2022 fund holdings, matching the layout of
#    the other quarterly sheets (2021-Q4, 2021-Q3, ...).
# 2. Prepend a new summary row for "2022-Q1" at the top of the "总计" sheet's
#    data (pushing the existing rows down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q1" worksheet, inserted before "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$ws = $wb.Worksheets.Add($totalSheet)
$ws.Name = "2022-Q1"

# NOTE: inserting a sheet "before" $totalSheet reseats that anchor onto the
# newly-created sheet (position-based anchor, not object identity) - so we
# must re-resolve "总计" by name before touching it again below.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy the header formatting (bold / centered / bordered) from an existing
# quarterly sheet so the new sheet matches the workbook's house style
# exactly (reuses the same cell style instead of inventing a new one).
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$templateSheet.Range("A1:H1").Copy()
$ws.Range("A1:H1").PasteSpecial(-4122)
$templateSheet.Range("A2:A2").Copy()
$ws.Range("A2:A10").PasteSpecial(-4122)

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# The source data keeps code / name / scale / position columns as plain
# text (so fund codes keep their leading zeros and values like "99.00"
# keep their trailing zero) - force Text format before writing them.
$ws.Range("B2:G10").NumberFormat = "@"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "009909"
$ws.Range("C2").Value = "嘉实动力先锋混合A"
$ws.Range("D2").Value = "29.43"
$ws.Range("E2").Value = "90.99"
$ws.Range("F2").Value = "4.07"
$ws.Range("G2").Value = "1.1978"
$ws.Range("H2").Value = 8

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "009869"
$ws.Range("C3").Value = "嘉实产业先锋混合A"
$ws.Range("D3").Value = "17.16"
$ws.Range("E3").Value = "89.91"
$ws.Range("F3").Value = "3.83"
$ws.Range("G3").Value = "0.6572"
$ws.Range("H3").Value = 9

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "009870"
$ws.Range("C4").Value = "嘉实产业先锋混合C"
$ws.Range("D4").Value = "3.71"
$ws.Range("E4").Value = "89.91"
$ws.Range("F4").Value = "3.83"
$ws.Range("G4").Value = "0.1421"
$ws.Range("H4").Value = 9

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "009910"
$ws.Range("C5").Value = "嘉实动力先锋混合C"
$ws.Range("D5").Value = "2.95"
$ws.Range("E5").Value = "90.99"
$ws.Range("F5").Value = "4.07"
$ws.Range("G5").Value = "0.1201"
$ws.Range("H5").Value = 8

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "159792"
$ws.Range("C6").Value = "富国中证港股通互联网ETF"
$ws.Range("D6").Value = "2.76"
$ws.Range("E6").Value = "99.00"
$ws.Range("F6").Value = "3.65"
$ws.Range("G6").Value = "0.1007"
$ws.Range("H6").Value = 7

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "007151"
$ws.Range("C7").Value = "前海开源沪港深聚瑞混合"
$ws.Range("D7").Value = "0.60"
$ws.Range("E7").Value = "72.90"
$ws.Range("F7").Value = "6.54"
$ws.Range("G7").Value = "0.0392"
$ws.Range("H7").Value = 5

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "003993"
$ws.Range("C8").Value = "前海开源沪港深核心驱动灵活配置混合"
$ws.Range("D8").Value = "0.58"
$ws.Range("E8").Value = "82.10"
$ws.Range("F8").Value = "2.95"
$ws.Range("G8").Value = "0.0171"
$ws.Range("H8").Value = 10

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "513160"
$ws.Range("C9").Value = "银华恒生港股通中国科技ETF"
$ws.Range("D9").Value = "0.62"
$ws.Range("E9").Value = "92.07"
$ws.Range("F9").Value = "2.66"
$ws.Range("G9").Value = "0.0165"
$ws.Range("H9").Value = 10

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "006477"
$ws.Range("C10").Value = "中邮沪港深精选混合"
$ws.Range("D10").Value = "0.05"
$ws.Range("E10").Value = "83.91"
$ws.Range("F10").Value = "4.42"
$ws.Range("G10").Value = "0.0022"
$ws.Range("H10").Value = 8

# ---------------------------------------------------------------------
# 2) Insert the "2022-Q1" summary row at the top of the "总计" sheet
# ---------------------------------------------------------------------
$totalSheet.Range("A2:D2").EntireRow.Insert()

# New row 2 inherits formatting from the row that was pushed down to row 3
# (i.e. the old row 2) so it matches the rest of the table exactly.
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 9
$totalSheet.Range("D2").Value = 2.29

# Renumber the leading index column (A) sequentially for the rows that got
# pushed down, same as the rest of the table.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4

Write-Host "2022-Q1 sheet added; 总计 updated"
